$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")

# Update the E3-H2 user story (rows 25-26) to the new "filter" wording.
$ws.Range("C25").Value2 = "usuario"
$ws.Range("D25").Value2 = "Filtro catalogo"
$ws.Range("E25").Value2 = "Quiere ver los libros organizados utilizando filtros que me permitan ver solo libros de una categoria o genero especifico"
$ws.Range("G25").Value2 = "Filtro categorias"
$ws.Range("H25").Value2 = "quiero ver libros pertenecientes a una categoria especifica"
$ws.Range("I25").Value2 = "El usuario selecciona una categoria"
$ws.Range("J25").Value2 = "El sistema muestra solo los libros de la categoria elegida"

$ws.Range("G26").Value2 = "Filtro generos"
$ws.Range("H26").Value2 = "Quiero ver libros pertenecientes a un genero en especifico"
$ws.Range("I26").Value2 = "El usuario selecciona un genero"
$ws.Range("J26").Value2 = "El sistema muestra solo los libros del genero elegido"

# Restore the rows' natural (non-custom) height after editing their text,
# since the rows carried no explicit row height before the edit either.
$ws.Rows.Item(25).AutoFit()
$ws.Rows.Item(26).AutoFit()

# Unhide the previously hidden user-story rows (4-39) so they show in the sheet.
for ($i = 4; $i -le 39; $i++) {
    $ws.Rows.Item($i).Hidden = $false
}
